$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new price values are plain decimal numbers need to be
# forced to Text format first, otherwise Excel would silently convert
# them to numbers and drop meaningful trailing zeros (e.g. "0.0500").
$ws.Range("D5,D9,D16,D18,D19,D23,D25,D27,D30,D33,D34,D36,D37,D38,D41,D43,D46,D47,D48,D51").NumberFormat = "@"

$ws.Range("D2").Value = "27.153.60"
$ws.Range("E2").Value = "  +0.77%  "
$ws.Range("D3").Value = "1.679.91"
$ws.Range("E3").Value = "  +0.37%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "215.38"
$ws.Range("E5").Value = "  +0.14%  "
$ws.Range("E6").Value = "  +0.29%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  +1.86%  "
$ws.Range("D9").Value = "21.39"
$ws.Range("E9").Value = "  +5.51%  "
$ws.Range("E10").Value = "  +0.67%  "
$ws.Range("D12").Value = "1.917.25"
$ws.Range("E12").Value = "  +0.39%  "
$ws.Range("D13").Value = "1.688.75"
$ws.Range("E13").Value = "  +0.65%  "
$ws.Range("E14").Value = "  +1.56%  "
$ws.Range("E15").Value = "  +1.96%  "
$ws.Range("D16").Value = "66.29"
$ws.Range("E16").Value = "  +1.00%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "27.137.82"
$ws.Range("E17").Value = "  +0.70%  "
$ws.Range("B18").Value = "BitcoinCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D18").Value = "239.60"
$ws.Range("E18").Value = "  +1.84%  "
$ws.Range("D19").Value = "8.11"
$ws.Range("E19").Value = "  -0.38%  "
$ws.Range("D20").Value = "0.0₃0746"
$ws.Range("E20").Value = "  +1.48%  "
$ws.Range("E21").Value = "  -0.02%  "
$ws.Range("E22").Value = "  +1.71%  "
$ws.Range("D23").Value = "9.48"
$ws.Range("E23").Value = "  +3.17%  "
$ws.Range("E24").Value = "  -2.69%  "
$ws.Range("D25").Value = "146.82"
$ws.Range("E25").Value = "  +0.80%  "
$ws.Range("E26").Value = "  +0.59%  "
$ws.Range("D27").Value = "16.36"
$ws.Range("E27").Value = "  +2.04%  "
$ws.Range("E28").Value = "  +0.67%  "
$ws.Range("E29").Value = "  +0.08%  "
$ws.Range("D30").Value = "0.0500"
$ws.Range("E30").Value = "  +0.43%  "
$ws.Range("E31").Value = "  -0.08%  "
$ws.Range("B32").Value = "Maker"
$ws.Range("C32").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D32").Value = "1.558.77"
$ws.Range("E32").Value = "  +4.69%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "3.37"
$ws.Range("E33").Value = "  +1.32%  "
$ws.Range("D34").Value = "3.22"
$ws.Range("E34").Value = "  +2.56%  "
$ws.Range("E35").Value = "  +0.77%  "
$ws.Range("D36").Value = "0.604"
$ws.Range("E36").Value = "  +3.59%  "
$ws.Range("B37").Value = "HuobiToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D37").Value = "2.39"
$ws.Range("E37").Value = "  -1.19%  "
$ws.Range("B38").Value = "ARBITRUM"
$ws.Range("C38").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D38").Value = "0.936"
$ws.Range("E38").Value = "  +4.45%  "
$ws.Range("E39").Value = "  +2.61%  "
$ws.Range("E40").Value = "  +0.77%  "
$ws.Range("D41").Value = "69.46"
$ws.Range("E41").Value = "  +2.97%  "
$ws.Range("E42").Value = "  -0.01%  "
$ws.Range("D43").Value = "5.67"
$ws.Range("E43").Value = "  -3.01%  "
$ws.Range("E44").Value = "  -2.63%  "
$ws.Range("D45").Value = "1.825.58"
$ws.Range("E45").Value = "  +0.60%  "
$ws.Range("D46").Value = "0.781"
$ws.Range("E46").Value = "  +0.56%  "
$ws.Range("D47").Value = "90.73"
$ws.Range("E47").Value = "  +0.20%  "
$ws.Range("D48").Value = "1.59"
$ws.Range("E48").Value = "  +3.53%  "
$ws.Range("D49").Value = "0.0₆0107"
$ws.Range("E49").Value = "  +1.40%  "
$ws.Range("E50").Value = "  +2.59%  "
$ws.Range("D51").Value = "8.12"
$ws.Range("E51").Value = "  +5.28%  "
